$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New summary statistics added in column G
$ws.Range("G3").Value = "Average samples done per user"
$ws.Range("G4").Formula = "=AVERAGE(A1:A17)"

$ws.Range("G6").Value = "Cumulative Mean"
$ws.Range("G7").Formula = "=D22"

$ws.Range("G9").Value = "Per Person std deviation"
$ws.Range("G10").Formula = "=STDEV.P(B1:B17)"

# Update the view: zoom to 87% and select G8 (this also moves the
# scrolled-to top-left cell back to A1, dropping topLeftCell="A7")
$ws.Application.ActiveWindow.Zoom = 87
$ws.Range("G8").Select()
